# Fix issue with die roll not showing up. Worked on initial combat calendar check.
#
# Updates the "e006 Combat Calendar Check" event text in cell B7 of the
# "Events" worksheet: separates the r4.1/Calendar buttons onto their own
# line, adds a short explanatory paragraph about rolling for combat and
# continuing to e007 / next day, reflows the "Date from Combat Calendar"
# line, and changes the probability comparison from "PROBABILITY > " to
# "PROBABILITY >= ". The added text makes the wrapped cell taller, so the
# row height is updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newText = "<Bold>e006 Combat Calendar Check</Bold> `n" +
    "<InlineUIContainer><Button Content='r4.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n" +
    "<InlineUIContainer><Button Content='Calendar' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>`n" +
    "<LineBreak/><LineBreak/>`n" +
    "Roll for possible combat today. If die &lt;= probability, start morning briefing per `n" +
    "<InlineUIContainer><Button Content='e007' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.  `n" +
    "Otherwise continue with next day check.`n" +
    "<LineBreak/><LineBreak/>`n" +
    "Date from Combat Calendar: DATE<LineBreak/>`n" +
    "Expected Resistance: RESISTANCE<LineBreak/>`n" +
    "Probablility of Combat: PROBABILITY &gt;= <InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer>`n" +
    "<LineBreak/><LineBreak/>`n"

$ws.Range("B7").Value = $newText

# The reflowed/expanded text needs a taller wrapped row (128.4 -> 185.45).
$ws.Rows.Item(7).RowHeight = 185.45
